# feat: added big onion and poison gas generation that hurts the player while inside
#
# Adds a new credits row (row 5) for the explosion sprite sheet used by the
# new poison-gas / explosion effect, matching the "File Name / Source /
# Licensing / Notes" table already on the sheet, and updates the current
# selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New credit row (row 5): explosion_01_strip13.png ----------------------
# Column order mirrors the header row: File Name | Source | Licensing | Notes/Other
$ws.Range("A5").Value = "explosion_01_strip13.png"
$ws.Range("D5").Value = '"Bleed - http://remusprites.carbonmade.com/"'
$ws.Range("C5").Value = "Creative Commons License 3"
$ws.Range("B5").Value = "https://opengameart.org/content/simple-explosion-bleeds-game-art"

# The "Notes/Other" cell (D5) uses the same 14pt Arial look as the other
# credit rows (e.g. C8) but with the theme (automatic) text colour rather
# than the dark-grey rgb one, so pick up that formatting first...
$ws.Range("C8").Copy()
$ws.Range("D5").PasteSpecial(-4122)  # xlPasteFormats
# ...then switch the font colour over to the theme colour.
$ws.Range("D5").Font.ThemeColor = 1  # xlThemeColorDark1

# Row 5 is a touch taller than the default, like the other credit rows.
$ws.Rows.Item(5).RowHeight = 18

# --- Selection -------------------------------------------------------------
$ws.Range("B8").Select() | Out-Null
